# Scheduled market-data refresh: update currentAveragePrice / LevePrice /
# LeveProfit figures across the Leve profit-tracking sheets (ALC, ARM, BSM,
# CRP, CUL, GSM, LTW) to reflect the latest Universalis price snapshot.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3670.25
$ws.Range("I69").Value = 2985.8333
$ws.Range("J69").Value = 4696.875
$ws.Range("K69").Value = 8957.499899999999
$ws.Range("L69").Value = 14090.625
$ws.Range("M69").Value = -8083.499899999999
$ws.Range("N69").Value = -15838.625
$ws.Range("H72").Value = 3670.25
$ws.Range("I72").Value = 2985.8333
$ws.Range("J72").Value = 4696.875
$ws.Range("K72").Value = 26872.4997
$ws.Range("L72").Value = 42271.875
$ws.Range("M72").Value = -22504.4997
$ws.Range("N72").Value = -51007.875
$ws.Range("H74").Value = 5136.364
$ws.Range("I74").Value = 4666.6665
$ws.Range("J74").Value = 5312.5
$ws.Range("K74").Value = 4666.6665
$ws.Range("L74").Value = 5312.5
$ws.Range("M74").Value = -3730.6665
$ws.Range("N74").Value = -7184.5
$ws.Range("H77").Value = 5136.364
$ws.Range("I77").Value = 4666.6665
$ws.Range("J77").Value = 5312.5
$ws.Range("K77").Value = 23333.3325
$ws.Range("L77").Value = 26562.5
$ws.Range("M77").Value = -18653.3325
$ws.Range("N77").Value = -35922.5
$ws.Range("H100").Value = 2972.4138
$ws.Range("I100").Value = 2566.6667
$ws.Range("J100").Value = 4037.5
$ws.Range("K100").Value = 2566.6667
$ws.Range("L100").Value = 4037.5
$ws.Range("M100").Value = -2025.6667
$ws.Range("N100").Value = -5119.5
$ws.Range("H107").Value = 453.21738
$ws.Range("I107").Value = 239.88235
$ws.Range("J107").Value = 1057.6666
$ws.Range("K107").Value = 239.88235
$ws.Range("L107").Value = 1057.6666
$ws.Range("M107").Value = 1680.11765
$ws.Range("N107").Value = -4897.6666
$ws.Range("H113").Value = 6791.472
$ws.Range("I113").Value = 2504.9473
$ws.Range("J113").Value = 11582.294
$ws.Range("K113").Value = 2504.9473
$ws.Range("L113").Value = 11582.294
$ws.Range("M113").Value = 749.0527000000002
$ws.Range("N113").Value = -18090.294
$ws.Range("H129").Value = 970.5
$ws.Range("J129").Value = 1074.8334
$ws.Range("L129").Value = 3224.5002
$ws.Range("N129").Value = -13224.5002

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1089.7778
$ws.Range("I2").Value = 976
$ws.Range("K2").Value = 976
$ws.Range("M2").Value = -863
$ws.Range("H110").Value = 50599.5
$ws.Range("I110").Value = 63072.75
$ws.Range("J110").Value = 706.5
$ws.Range("K110").Value = 63072.75
$ws.Range("L110").Value = 706.5
$ws.Range("M110").Value = -61027.75
$ws.Range("N110").Value = -4796.5
$ws.Range("H116").Value = 1089.7778
$ws.Range("I116").Value = 976
$ws.Range("K116").Value = 976
$ws.Range("M116").Value = 1318

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1089.7778
$ws.Range("I3").Value = 976
$ws.Range("K3").Value = 976
$ws.Range("M3").Value = -862

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1632.875
$ws.Range("I16").Value = 1692.8572
$ws.Range("K16").Value = 1692.8572
$ws.Range("M16").Value = -1405.8572
$ws.Range("H113").Value = 1632.875
$ws.Range("I113").Value = 1692.8572
$ws.Range("K113").Value = 1692.8572
$ws.Range("M113").Value = 477.1428000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 4132.75
$ws.Range("I81").Value = 750.5
$ws.Range("K81").Value = 2251.5
$ws.Range("M81").Value = -1128.5
$ws.Range("H84").Value = 4132.75
$ws.Range("I84").Value = 750.5
$ws.Range("K84").Value = 6754.5
$ws.Range("M84").Value = -1138.5
$ws.Range("H97").Value = 1280.5
$ws.Range("J97").Value = 1250.75
$ws.Range("L97").Value = 3752.25
$ws.Range("N97").Value = -4744.25
$ws.Range("H100").Value = 5139.9
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 5305.1577
$ws.Range("K100").Value = 6000
$ws.Range("L100").Value = 15915.4731
$ws.Range("M100").Value = -5189
$ws.Range("N100").Value = -17537.4731
$ws.Range("H113").Value = 789.6667
$ws.Range("I113").Value = 547.6
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1642.8
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = 527.1999999999998
$ws.Range("N113").Value = -10340
$ws.Range("H122").Value = 900.9143
$ws.Range("I122").Value = 501.22223
$ws.Range("K122").Value = 4511.00007
$ws.Range("M122").Value = -2061.00007
$ws.Range("H131").Value = 853.25
$ws.Range("I131").Value = 288
$ws.Range("K131").Value = 864
$ws.Range("M131").Value = 4176

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2749.375
$ws.Range("I80").Value = 2832.5
$ws.Range("K80").Value = 2832.5
$ws.Range("M80").Value = -1834.5
$ws.Range("H82").Value = 39794
$ws.Range("J82").Value = 39794
$ws.Range("L82").Value = 39794
$ws.Range("N82").Value = -40560
$ws.Range("H83").Value = 2749.375
$ws.Range("I83").Value = 2832.5
$ws.Range("K83").Value = 14162.5
$ws.Range("M83").Value = -9170.5
$ws.Range("H85").Value = 39794
$ws.Range("J85").Value = 39794
$ws.Range("L85").Value = 39794
$ws.Range("N85").Value = -42446
$ws.Range("H113").Value = 1820.1
$ws.Range("I113").Value = 1628.5714
$ws.Range("J113").Value = 2267
$ws.Range("K113").Value = 1628.5714
$ws.Range("L113").Value = 2267
$ws.Range("M113").Value = 541.4286
$ws.Range("N113").Value = -6607

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2208.2
$ws.Range("I61").Value = 2616.4
$ws.Range("K61").Value = 2616.4
$ws.Range("M61").Value = -2414.4
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = $null
$ws.Range("N81").Value = $null
$ws.Range("H82").Value = 1392.5
$ws.Range("I82").Value = 1392.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1392.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1031.5
$ws.Range("N82").Value = $null
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = $null
$ws.Range("N84").Value = $null
$ws.Range("H85").Value = 1392.5
$ws.Range("I85").Value = 1392.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1392.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -144.5
$ws.Range("N85").Value = $null
$ws.Range("H113").Value = 2208.2
$ws.Range("I113").Value = 2616.4
$ws.Range("K113").Value = 2616.4
$ws.Range("M113").Value = -446.4000000000001
$ws.Range("H136").Value = 15154553
$ws.Range("I136").Value = 3058.5
$ws.Range("J136").Value = 55558540
$ws.Range("K136").Value = 9175.5
$ws.Range("L136").Value = 166675620
$ws.Range("M136").Value = -6625.5
$ws.Range("N136").Value = -166680720
